# Weekly fruit/vegetable price update: a new observation row is inserted
# at row 54 (pushing the existing rows 54-72 down to 55-73), and the sheet's
# used range grows from A1:R72 to A1:R73.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54; Excel shifts rows 54:72 down to 55:73
# and the new row inherits formatting (e.g. the date style) from the row above.
$ws.Rows(54).Insert()

# Populate the newly inserted row 54 with the new market observation.
$ws.Cells.Item(54, 1).Value  = 10
$ws.Cells.Item(54, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(54, 3).Value  = "La Araucanía"
$ws.Cells.Item(54, 4).Value  = 44845
$ws.Cells.Item(54, 5).Value  = 9
$ws.Cells.Item(54, 6).Value  = 300000000
$ws.Cells.Item(54, 7).Value  = "Espárragos"
$ws.Cells.Item(54, 8).Value  = "Sin especificar"
$ws.Cells.Item(54, 9).Value  = "Primera"
$ws.Cells.Item(54, 10).Value = 110
$ws.Cells.Item(54, 11).Value = 1500
$ws.Cells.Item(54, 12).Value = 1600
$ws.Cells.Item(54, 13).Value = 1545
$ws.Cells.Item(54, 14).Value = "$/kilo"
$ws.Cells.Item(54, 15).Value = "Región del Maule"
$ws.Cells.Item(54, 16).Value = 1545
$ws.Cells.Item(54, 17).Value = 1
$ws.Cells.Item(54, 18).Value = "Hortaliza"
